$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Add the new "Interest:" heading paragraph at the end of the body.
# ------------------------------------------------------------------
$headingPara = $d.Paragraphs.Add()
$headingPara.Range.Text = "Interest:"
$headingPara.Range.ParagraphFormat.Style = "Heading1"

# ------------------------------------------------------------------
# 2. Add the body paragraph that follows it. We temporarily append a
#    unique marker after the real text so we can compute the exact
#    "end of real text" offset (the COM layer mis-handles bookmarks
#    whose target range collapses onto the very last character of the
#    document, so we give it some breathing room and trim the marker
#    off afterwards).
# ------------------------------------------------------------------
$bodyPara = $d.Paragraphs.Add()
$bodyPara.Range.ParagraphFormat.Style = "Normal"
$bodyText = "This would interest anyone who has a similar need and is in search of a rent or relocation. He will get a clear picture of his neighbourhood. By changing the locations or preferences, this algorithm can be used by anyone."
$marker = "ZZZ_TEMP_MARKER_ZZZ"
$bodyPara.Range.Text = $bodyText + $marker

# Position right after the real text (i.e. right before the marker,
# and before the paragraph mark).
$endOfRealText = $bodyPara.Range.End - 1 - $marker.Length
$target = $d.Range($endOfRealText, $endOfRealText)

# ------------------------------------------------------------------
# 3. Relocate the "_GoBack" bookmark to that collapsed range. Adding a
#    bookmark with a name that already exists moves it (Word keeps
#    bookmark names unique), so this both removes it from the old
#    "Business Problem:" heading and places it at the end of the new
#    paragraph's text.
# ------------------------------------------------------------------
$d.Bookmarks.Add("_GoBack", $target)

# ------------------------------------------------------------------
# 4. Strip the temporary marker text back out.
# ------------------------------------------------------------------
$null = $d.Content.Find.Execute($marker, $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
